# Auto-generated edit script applying the target diff to before.xlsx.
# The edit re-sorts/re-numbers the observation records: each data row (2,3,4,
# 6,7,8,9,10) ends up holding the field values that another row held before
# (row 5 is untouched). Rather than physically moving ranges (which risks
# clobbering data while rows are still being read), every destination cell is
# written directly with its final literal value taken from the target state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111486385
$ws.Range("I2").Value = "'20"
$ws.Range("Q2").Value = 624029.7289273632
$ws.Range("R2").Value = 6932998.597210908
$ws.Range("S2").Value = 10
$ws.Range("Z2").Value = "14:11"
$ws.Range("AB2").Value = "14:11"
$ws.Range("AC2").Value = "Ca 20 ex varav 3 blommande"

# Row 3
$ws.Range("A3").Value = 111486400
$ws.Range("I3").Value = "'5"
$ws.Range("K3").Value = "fullt utvecklade blad"
$ws.Range("P3").Value = "Glödenhöjden (Glödenhöjden), Mpd"
$ws.Range("Q3").Value = 624030.1824148977
$ws.Range("R3").Value = 6932961.620511409
$ws.Range("S3").Value = 10
$ws.Range("Y3").Value = "'2023-08-14"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").Value = "'2023-08-14"
$ws.Range("AB3").Value = "00:00"

# Row 4
$ws.Range("A4").Value = 111486117
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = "'10"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("K4").Value = "blomning"
$ws.Range("N4").Value = "observerad"
$ws.Range("P4").Value = "Glödenhöjden (Glödenhöjden), Mpd"
$ws.Range("Q4").Value = 623993.6707231236
$ws.Range("R4").Value = 6933021.760048959
$ws.Range("S4").Value = 15
$ws.Range("Z4").Value = "14:17"
$ws.Range("AB4").Value = "14:17"
$ws.Range("AC4").Value = "10 plantor varav 2 blommande"
$ws.Range("AH4").ClearContents() | Out-Null

# Row 6
$ws.Range("A6").Value = 111486415
$ws.Range("I6").Value = "'10"
$ws.Range("K6").Value = "fullt utvecklade blad"
$ws.Range("Q6").Value = 624040.2038791699
$ws.Range("R6").Value = 6932953.67081845
$ws.Range("Z6").Value = "13:46"
$ws.Range("AB6").Value = "13:46"
$ws.Range("AC6").ClearContents() | Out-Null

# Row 7
$ws.Range("A7").Value = 111485917
$ws.Range("I7").Value = "'3"
$ws.Range("K7").Value = "blomning"
$ws.Range("P7").Value = "Glödenhöjden nordost (Glödenhöjden), Mpd"
$ws.Range("Q7").Value = 624090.1097011974
$ws.Range("R7").Value = 6933043.392863069
$ws.Range("S7").Value = 25
$ws.Range("Z7").Value = "15:00"
$ws.Range("AB7").Value = "15:00"

# Row 8
$ws.Range("A8").Value = 111486450
$ws.Range("I8").Value = "'2"
$ws.Range("K8").Value = "fullt utvecklade blad"
$ws.Range("Q8").Value = 624051.1502826829
$ws.Range("R8").Value = 6932945.755648845
$ws.Range("Z8").Value = "13:43"
$ws.Range("AB8").Value = "13:43"

# Row 9
$ws.Range("A9").Value = 111485854
$ws.Range("B9").Value = 89405
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I9").Value = "'1"
$ws.Range("J9").Value = "fruktkroppar"
$ws.Range("K9").ClearContents() | Out-Null
$ws.Range("N9").ClearContents() | Out-Null
$ws.Range("P9").Value = "Glödenhöjden nordost (Glödenhöjden), Mpd"
$ws.Range("Q9").Value = 624096.1730324102
$ws.Range("R9").Value = 6933042.231978768
$ws.Range("S9").Value = 20
$ws.Range("Z9").Value = "15:02"
$ws.Range("AB9").Value = "15:02"
$ws.Range("AH9").Value = "Häll- och blockskog"

# Row 10
$ws.Range("A10").Value = 111486280
$ws.Range("I10").Value = "'3"
$ws.Range("J10").ClearContents() | Out-Null
$ws.Range("K10").Value = "blomning"
$ws.Range("Q10").Value = 624009.7035872869
$ws.Range("R10").Value = 6933014.034667149
$ws.Range("Y10").Value = "'2023-08-12"
$ws.Range("Z10").Value = "14:17"
$ws.Range("AA10").Value = "'2023-08-12"
$ws.Range("AB10").Value = "14:17"
